$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "甲數/乙數" sentence-building table in A8:C12 loses its first data
# row (which held the stray "數" fragment): every remaining row's content
# slides up by one, and the final row of each column is cleared out.
#   A8:C8  <- old A9:C9      (多/bigger/大)
#   A9:C9  <- old A10:C10    (5/int/數量)
#   A10:B10 <- old A11:B11   (。/x)
#   A11, C10, B12 end up empty
$ws.Range("A9:C9").Copy($ws.Range("A8:C8"))
$ws.Range("A10:C10").Copy($ws.Range("A9:C9"))
$ws.Range("A11:B11").Copy($ws.Range("A10:B10"))

$ws.Range("A11").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("B12").ClearContents()

# Reflect the new editing focus in the saved selection.
[void]$ws.Range("B8:C11").Select()
